$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force literal-text entry for values that would
# otherwise be auto-parsed as numbers by Excel (e.g. "1.000", "235.46").
# Written as Text, then copied + PasteSpecial(values-only) into the target
# cell so the target keeps its original (default) style/number format.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

# Row 2
$ws.Cells.Item(2,4).Value = "30.373.27"
$ws.Cells.Item(2,5).Value = "  +0.25%  "

# Row 3
$ws.Cells.Item(3,4).Value = "1.872.33"
$ws.Cells.Item(3,5).Value = "  -0.56%  "

# Row 4
$scratch.Value = "1.000"
$scratch.Copy() | Out-Null
$ws.Cells.Item(4,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(4,5).Value = "  -0.07%  "

# Row 5
$scratch.Value = "235.46"
$scratch.Copy() | Out-Null
$ws.Cells.Item(5,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(5,5).Value = "  -1.07%  "

# Row 6
$scratch.Value = "1.000"
$scratch.Copy() | Out-Null
$ws.Cells.Item(6,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(6,5).Value = "  -0.09%  "

# Row 7
$scratch.Value = "0.4673"
$scratch.Copy() | Out-Null
$ws.Cells.Item(7,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(7,5).Value = "  -0.02%  "

# Row 8
$scratch.Value = "0.2847"
$scratch.Copy() | Out-Null
$ws.Cells.Item(8,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(8,5).Value = "  +1.07%  "

# Row 9
$scratch.Value = "0.06572"
$scratch.Copy() | Out-Null
$ws.Cells.Item(9,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(9,5).Value = "  +0.31%  "

# Row 10
$scratch.Value = "21.48"
$scratch.Copy() | Out-Null
$ws.Cells.Item(10,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(10,5).Value = "  +9.19%  "

# Row 11
$scratch.Value = "0.07907"
$scratch.Copy() | Out-Null
$ws.Cells.Item(11,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(11,5).Value = "  +1.87%  "

# Row 12
$scratch.Value = "97.81"
$scratch.Copy() | Out-Null
$ws.Cells.Item(12,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(12,5).Value = "  +0.00%  "

# Row 13
$ws.Cells.Item(13,4).Value = "1.866.41"
$ws.Cells.Item(13,5).Value = "  -1.00%  "

# Row 14
$scratch.Value = "5.129"
$scratch.Copy() | Out-Null
$ws.Cells.Item(14,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(14,5).Value = "  +0.54%  "

# Row 15
$scratch.Value = "0.6772"
$scratch.Copy() | Out-Null
$ws.Cells.Item(15,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(15,5).Value = "  +1.95%  "

# Row 16
$scratch.Value = "281.37"
$scratch.Copy() | Out-Null
$ws.Cells.Item(16,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(16,5).Value = "  -0.42%  "

# Row 17
$ws.Cells.Item(17,4).Value = "30.375.21"
$ws.Cells.Item(17,5).Value = "  +0.24%  "

# Row 18
$scratch.Value = "0.9997"
$scratch.Copy() | Out-Null
$ws.Cells.Item(18,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(18,5).Value = "  -0.04%  "

# Row 19
$scratch.Value = "5.513"
$scratch.Copy() | Out-Null
$ws.Cells.Item(19,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(19,5).Value = "  +3.20%  "

# Row 20
$scratch.Value = "12.70"
$scratch.Copy() | Out-Null
$ws.Cells.Item(20,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(20,5).Value = "  +0.93%  "

# Row 21
$ws.Cells.Item(21,4).Value = "2.115.38"
$ws.Cells.Item(21,5).Value = "  -0.77%  "

# Row 22
$scratch.Value = "0.000007311"
$scratch.Copy() | Out-Null
$ws.Cells.Item(22,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(22,5).Value = "  +0.29%  "

# Row 23
$scratch.Value = "1.000"
$scratch.Copy() | Out-Null
$ws.Cells.Item(23,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(23,5).Value = "  -0.19%  "

# Row 24
$scratch.Value = "6.212"
$scratch.Copy() | Out-Null
$ws.Cells.Item(24,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(24,5).Value = "  +0.89%  "

# Row 25
$scratch.Value = "9.287"
$scratch.Copy() | Out-Null
$ws.Cells.Item(25,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(25,5).Value = "  +0.26%  "

# Row 26
$scratch.Value = "165.07"
$scratch.Copy() | Out-Null
$ws.Cells.Item(26,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(26,5).Value = "  -1.25%  "

# Row 27
$scratch.Value = "19.17"
$scratch.Copy() | Out-Null
$ws.Cells.Item(27,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(27,5).Value = "  +0.73%  "

# Row 28
$scratch.Value = "1.941"
$scratch.Copy() | Out-Null
$ws.Cells.Item(28,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(28,5).Value = "  -2.23%  "

# Row 29
$ws.Cells.Item(29,5).Value = "  +0.28%  "

# Row 30
$scratch.Value = "0.09712"
$scratch.Copy() | Out-Null
$ws.Cells.Item(30,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(30,5).Value = "  -0.80%  "

# Row 31
$scratch.Value = "4.432"
$scratch.Copy() | Out-Null
$ws.Cells.Item(31,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(31,5).Value = "  -0.42%  "

# Row 32
$ws.Cells.Item(32,5).Value = "  -0.97%  "

# Row 33
$scratch.Value = "4.119"
$scratch.Copy() | Out-Null
$ws.Cells.Item(33,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(33,5).Value = "  -1.00%  "

# Row 34
$scratch.Value = "0.04695"
$scratch.Copy() | Out-Null
$ws.Cells.Item(34,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(34,5).Value = "  +0.06%  "

# Row 35
$scratch.Value = "1.119"
$scratch.Copy() | Out-Null
$ws.Cells.Item(35,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(35,5).Value = "  +2.61%  "

# Row 36
$scratch.Value = "0.7080"
$scratch.Copy() | Out-Null
$ws.Cells.Item(36,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(36,5).Value = "  +0.23%  "

# Row 37
$scratch.Value = "2.711"
$scratch.Copy() | Out-Null
$ws.Cells.Item(37,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(37,5).Value = "  -0.25%  "

# Row 38
$scratch.Value = "0.01865"
$scratch.Copy() | Out-Null
$ws.Cells.Item(38,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(38,5).Value = "  -0.04%  "

# Row 39
$scratch.Value = "6.330"
$scratch.Copy() | Out-Null
$ws.Cells.Item(39,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(39,5).Value = "  -5.89%  "

# Row 40
$scratch.Value = "2.541"
$scratch.Copy() | Out-Null
$ws.Cells.Item(40,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(40,5).Value = "  +0.75%  "

# Row 41
$scratch.Value = "73.25"
$scratch.Copy() | Out-Null
$ws.Cells.Item(41,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(41,5).Value = "  +1.70%  "

# Row 42
$scratch.Value = "1.948"
$scratch.Copy() | Out-Null
$ws.Cells.Item(42,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(42,5).Value = "  -0.88%  "

# Row 43
$scratch.Value = "0.8500"
$scratch.Copy() | Out-Null
$ws.Cells.Item(43,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(43,5).Value = "  -2.44%  "

# Row 44
$scratch.Value = "0.4196"
$scratch.Copy() | Out-Null
$ws.Cells.Item(44,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(44,5).Value = "  +0.19%  "

# Row 45
$ws.Cells.Item(45,2).Value = "Quant"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$scratch.Value = "103.95"
$scratch.Copy() | Out-Null
$ws.Cells.Item(45,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(45,5).Value = "  -0.07%  "

# Row 46
$ws.Cells.Item(46,2).Value = "PaxDollar"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$scratch.Value = "1.000"
$scratch.Copy() | Out-Null
$ws.Cells.Item(46,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(46,5).Value = "  -0.07%  "

# Row 47
$scratch.Value = "7.225"
$scratch.Copy() | Out-Null
$ws.Cells.Item(47,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(47,5).Value = "  +0.27%  "

# Row 48
$scratch.Value = "9.147"
$scratch.Copy() | Out-Null
$ws.Cells.Item(48,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(48,5).Value = "  -1.23%  "

# Row 49
$scratch.Value = "932.85"
$scratch.Copy() | Out-Null
$ws.Cells.Item(49,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(49,5).Value = "  -5.35%  "

# Row 50
$scratch.Value = "34.14"
$scratch.Copy() | Out-Null
$ws.Cells.Item(50,4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(50,5).Value = "  +0.39%  "

# Row 51
$ws.Cells.Item(51,5).Value = "  -2.35%  "

$scratch.Clear() | Out-Null
Write-Output "done"